$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.660594940185547
$ws.Range("B1").Value = 3.721952676773071
$ws.Range("C1").Value = 2.881871461868286
$ws.Range("D1").Value = 2.756474733352661
$ws.Range("E1").Value = 2.52364706993103
